# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 721.7143
$ws.Range("I28").Value = 522.7778
$ws.Range("K28").Value = 522.7778
$ws.Range("M28").Value = -37.77779999999996
$ws.Range("H29").Value = 8497.666999999999
$ws.Range("J29").Value = 10497
$ws.Range("L29").Value = 31491
$ws.Range("N29").Value = -32053
$ws.Range("H32").Value = 6215.4546
$ws.Range("I32").Value = 3472.75
$ws.Range("J32").Value = 7782.7144
$ws.Range("K32").Value = 3472.75
$ws.Range("L32").Value = 7782.7144
$ws.Range("M32").Value = -3146.75
$ws.Range("N32").Value = -8434.714400000001
$ws.Range("H69").Value = 11714.214
$ws.Range("J69").Value = 12333.333
$ws.Range("L69").Value = 36999.999
$ws.Range("N69").Value = -38747.999
$ws.Range("H72").Value = 11714.214
$ws.Range("J72").Value = 12333.333
$ws.Range("L72").Value = 110999.997
$ws.Range("N72").Value = -119735.997
$ws.Range("H94").Value = 919.5
$ws.Range("I94").Value = 919.5
$ws.Range("K94").Value = 919.5
$ws.Range("M94").Value = -468.5
$ws.Range("H107").Value = 1107.5
$ws.Range("I107").Value = 1373.5
$ws.Range("K107").Value = 1373.5
$ws.Range("M107").Value = 546.5
$ws.Range("H113").Value = 52499.75
$ws.Range("I113").Value = 52500
$ws.Range("J113").Value = 52499.5
$ws.Range("K113").Value = 52500
$ws.Range("L113").Value = 52499.5
$ws.Range("M113").Value = -49246
$ws.Range("N113").Value = -59007.5
$ws.Range("H115").Value = 981.5
$ws.Range("I115").Value = 981.5
$ws.Range("K115").Value = 2944.5
$ws.Range("M115").Value = -1377.5
$ws.Range("H132").Value = 4464.15
$ws.Range("I132").Value = 4172.816
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 12518.448
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -9988.448
$ws.Range("N132").Value = -35058.5
$ws.Range("H135").Value = 1846.92
$ws.Range("I135").Value = 1733.6522
$ws.Range("J135").Value = 3149.5
$ws.Range("K135").Value = 15602.8698
$ws.Range("L135").Value = 28345.5
$ws.Range("M135").Value = -13067.8698
$ws.Range("N135").Value = -33415.5
$ws.Range("H137").Value = 2969.4866
$ws.Range("I137").Value = 2803.1667
$ws.Range("K137").Value = 8409.500100000001
$ws.Range("M137").Value = -5859.500100000001
$ws.Range("H138").Value = 3524.7
$ws.Range("J138").Value = 3475.3022
$ws.Range("L138").Value = 10425.9066
$ws.Range("N138").Value = -20705.9066

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12066.5
$ws.Range("I32").Value = 9405.166999999999
$ws.Range("K32").Value = 9405.166999999999
$ws.Range("M32").Value = -9118.166999999999
$ws.Range("H97").Value = 893.13336
$ws.Range("I97").Value = 869.0769
$ws.Range("J97").Value = 1049.5
$ws.Range("K97").Value = 869.0769
$ws.Range("L97").Value = 1049.5
$ws.Range("M97").Value = -373.0769
$ws.Range("N97").Value = -2041.5
$ws.Range("H122").Value = 12149.8
$ws.Range("I122").Value = 14545.272
$ws.Range("K122").Value = 43635.81600000001
$ws.Range("M122").Value = -41185.81600000001
$ws.Range("H132").Value = 4331.1577
$ws.Range("I132").Value = 3356.2
$ws.Range("K132").Value = 10068.6
$ws.Range("M132").Value = -7538.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 11900
$ws.Range("J33").Value = 23300
$ws.Range("L33").Value = 23300
$ws.Range("N33").Value = -23972
$ws.Range("H60").Value = 36994.75
$ws.Range("J60").Value = 36994.75
$ws.Range("L60").Value = 36994.75
$ws.Range("N60").Value = -38192.75
$ws.Range("H134").Value = 3780.3333
$ws.Range("I134").Value = 2920.75
$ws.Range("J134").Value = 5499.5
$ws.Range("K134").Value = 8762.25
$ws.Range("L134").Value = 16498.5
$ws.Range("M134").Value = -6227.25
$ws.Range("N134").Value = -21568.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 52323.92
$ws.Range("I16").Value = 53021.2
$ws.Range("J16").Value = 49999.668
$ws.Range("K16").Value = 53021.2
$ws.Range("L16").Value = 49999.668
$ws.Range("M16").Value = -52734.2
$ws.Range("N16").Value = -50573.668
$ws.Range("H31").Value = 6548.1035
$ws.Range("I31").Value = 6290.25
$ws.Range("J31").Value = 6970.0454
$ws.Range("K31").Value = 6290.25
$ws.Range("L31").Value = 6970.0454
$ws.Range("M31").Value = -5995.25
$ws.Range("N31").Value = -7560.0454
$ws.Range("H34").Value = 6548.1035
$ws.Range("I34").Value = 6290.25
$ws.Range("J34").Value = 6970.0454
$ws.Range("K34").Value = 6290.25
$ws.Range("L34").Value = 6970.0454
$ws.Range("M34").Value = -6088.25
$ws.Range("N34").Value = -7374.0454
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("H58").Value = 6237.6
$ws.Range("I58").Value = 6332.1577
$ws.Range("K58").Value = 6332.1577
$ws.Range("M58").Value = -6129.1577
$ws.Range("H113").Value = 52323.92
$ws.Range("I113").Value = 53021.2
$ws.Range("J113").Value = 49999.668
$ws.Range("K113").Value = 53021.2
$ws.Range("L113").Value = 49999.668
$ws.Range("M113").Value = -50851.2
$ws.Range("N113").Value = -54339.668
$ws.Range("H134").Value = 2258.0908
$ws.Range("I134").Value = 2033.9
$ws.Range("K134").Value = 6101.700000000001
$ws.Range("M134").Value = -3566.700000000001
$ws.Range("H136").Value = 6237.6
$ws.Range("I136").Value = 6332.1577
$ws.Range("K136").Value = 18996.4731
$ws.Range("M136").Value = -16446.4731
$ws.Range("H141").Value = 159790.75
$ws.Range("J141").Value = 164564.27
$ws.Range("L141").Value = 164564.27
$ws.Range("N141").Value = -174924.27
$ws.Range("M56").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 10000
$ws.Range("J100").Value = 10000
$ws.Range("L100").Value = 30000
$ws.Range("N100").Value = -31622
$ws.Range("H103").Value = 596.3333
$ws.Range("I103").Value = 1075
$ws.Range("K103").Value = 3225
$ws.Range("M103").Value = -2346
$ws.Range("H113").Value = 3684.76
$ws.Range("I113").Value = 2067
$ws.Range("J113").Value = 4089.2
$ws.Range("K113").Value = 6201
$ws.Range("L113").Value = 12267.6
$ws.Range("M113").Value = -4031
$ws.Range("N113").Value = -16607.6
$ws.Range("H132").Value = 2874.5833
$ws.Range("J132").Value = 2928.5715
$ws.Range("L132").Value = 26357.1435
$ws.Range("N132").Value = -31417.1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 684.8333
$ws.Range("I97").Value = 687.619
$ws.Range("J97").Value = 665.3333
$ws.Range("K97").Value = 687.619
$ws.Range("L97").Value = 665.3333
$ws.Range("M97").Value = -191.619
$ws.Range("N97").Value = -1657.3333
$ws.Range("H113").Value = 25638.25
$ws.Range("I113").Value = 18999.334
$ws.Range("K113").Value = 18999.334
$ws.Range("M113").Value = -16829.334
$ws.Range("H126").Value = 4452.846
$ws.Range("J126").Value = 5316
$ws.Range("L126").Value = 15948
$ws.Range("N126").Value = -20888
$ws.Range("H132").Value = 6881.6665
$ws.Range("I132").Value = 6289.1724
$ws.Range("J132").Value = 8599.9
$ws.Range("K132").Value = 18867.5172
$ws.Range("L132").Value = 25799.7
$ws.Range("M132").Value = -16337.5172
$ws.Range("N132").Value = -30859.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4334.3335
$ws.Range("I7").Value = 3332.6667
$ws.Range("K7").Value = 3332.6667
$ws.Range("M7").Value = -3220.6667
$ws.Range("H22").Value = 4999
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4999
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4999
$ws.Range("N22").Value = -5589
$ws.Range("H27").Value = 4999
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4999
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4999
$ws.Range("N27").Value = -5213
$ws.Range("H126").Value = 4334.3335
$ws.Range("I126").Value = 3332.6667
$ws.Range("K126").Value = 9998.000100000001
$ws.Range("M126").Value = -7528.000100000001
$ws.Range("H132").Value = 2961.5
$ws.Range("I132").Value = 2840.0527
$ws.Range("K132").Value = 8520.158100000001
$ws.Range("M132").Value = -5990.158100000001
$ws.Range("M22").ClearContents()
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 199815
$ws.Range("J5").Value = 199815
$ws.Range("L5").Value = 199815
$ws.Range("N5").Value = -200039
$ws.Range("H100").Value = 2178.9656
$ws.Range("I100").Value = 2259.44
$ws.Range("K100").Value = 4518.88
$ws.Range("M100").Value = -3977.88
$ws.Range("H122").Value = 9650.166999999999
$ws.Range("I122").Value = 9298.666999999999
$ws.Range("J122").Value = 10001.667
$ws.Range("K122").Value = 27896.001
$ws.Range("L122").Value = 30005.001
$ws.Range("M122").Value = -25446.001
$ws.Range("N122").Value = -34905.001
$ws.Range("H126").Value = 1738.8064
$ws.Range("I126").Value = 1737.9656
$ws.Range("K126").Value = 5213.8968
$ws.Range("M126").Value = -2743.8968
$ws.Range("H132").Value = 7991.5835
$ws.Range("J132").Value = 14949.5
$ws.Range("L132").Value = 44848.5
$ws.Range("N132").Value = -49908.5
$ws.Range("H136").Value = 17866.379
$ws.Range("J136").Value = 3024.6843
$ws.Range("L136").Value = 9074.052899999999
$ws.Range("N136").Value = -14174.0529
